$wb = $excel.ActiveWorkbook

# Sheet 1 (ALC), row 19
$ws = $wb.Worksheets.Item(1)
$ws.Range("H19").Value = 1701.5
$ws.Range("I19").Value = 2247.6365
$ws.Range("J19").Value = 500
$ws.Range("K19").Value = 2247.6365
$ws.Range("L19").Value = 500
$ws.Range("M19").Value = -2072.6365
$ws.Range("N19").Value = -850

# Sheet 1 (ALC), row 87
$ws = $wb.Worksheets.Item(1)
$ws.Range("H87").Value = 47965.69
$ws.Range("J87").Value = 47965.69
$ws.Range("L87").Value = 47965.69
$ws.Range("N87").Value = -50461.69

# Sheet 1 (ALC), row 90
$ws = $wb.Worksheets.Item(1)
$ws.Range("H90").Value = 47965.69
$ws.Range("J90").Value = 47965.69
$ws.Range("L90").Value = 143897.07
$ws.Range("N90").Value = -156377.07

# Sheet 1 (ALC), row 107
$ws = $wb.Worksheets.Item(1)
$ws.Range("H107").Value = 1478.24
$ws.Range("I107").Value = 874.1429000000001
$ws.Range("J107").Value = 4649.75
$ws.Range("K107").Value = 874.1429000000001
$ws.Range("L107").Value = 4649.75
$ws.Range("M107").Value = 1045.8571
$ws.Range("N107").Value = -8489.75

# Sheet 1 (ALC), row 132
$ws = $wb.Worksheets.Item(1)
$ws.Range("H132").Value = 1873.9375
$ws.Range("I132").Value = 1427.5714
$ws.Range("K132").Value = 4282.7142
$ws.Range("M132").Value = -1752.7142

# Sheet 1 (ALC), row 138
$ws = $wb.Worksheets.Item(1)
$ws.Range("H138").Value = 4020.6057
$ws.Range("I138").Value = 2634.8667
$ws.Range("J138").Value = 5034.561
$ws.Range("K138").Value = 7904.6001
$ws.Range("L138").Value = 15103.683
$ws.Range("M138").Value = -2764.6001
$ws.Range("N138").Value = -25383.683

# Sheet 1 (ALC), row 141
$ws = $wb.Worksheets.Item(1)
$ws.Range("H141").Value = 2429.4375
$ws.Range("I141").Value = 2565.0715
$ws.Range("K141").Value = 7695.2145
$ws.Range("M141").Value = -2515.2145

# Sheet 2 (ARM), row 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = 4248.9644
$ws.Range("I2").Value = 4086.2917
$ws.Range("K2").Value = 4086.2917
$ws.Range("M2").Value = -3973.2917

# Sheet 2 (ARM), row 37
$ws = $wb.Worksheets.Item(2)
$ws.Range("H37").Value = 18258.25
$ws.Range("J37").Value = 18258.25
$ws.Range("L37").Value = 18258.25
$ws.Range("N37").Value = -18804.25

# Sheet 2 (ARM), row 55
$ws = $wb.Worksheets.Item(2)
$ws.Range("H55").Value = 40465.867
$ws.Range("J55").Value = 48082.332
$ws.Range("L55").Value = 48082.332
$ws.Range("N55").Value = -48712.332

# Sheet 2 (ARM), row 64
$ws = $wb.Worksheets.Item(2)
$ws.Range("H64").Value = 30091
$ws.Range("J64").Value = 30091
$ws.Range("L64").Value = 30091
$ws.Range("N64").Value = -30587

# Sheet 2 (ARM), row 67
$ws = $wb.Worksheets.Item(2)
$ws.Range("H67").Value = 30091
$ws.Range("J67").Value = 30091
$ws.Range("L67").Value = 30091
$ws.Range("N67").Value = -31807

# Sheet 2 (ARM), row 80
$ws = $wb.Worksheets.Item(2)
$ws.Range("H80").Value = 41166.668
$ws.Range("J80").Value = 41166.668
$ws.Range("L80").Value = 41166.668
$ws.Range("N80").Value = -43162.668

# Sheet 2 (ARM), row 83
$ws = $wb.Worksheets.Item(2)
$ws.Range("H83").Value = 41166.668
$ws.Range("J83").Value = 41166.668
$ws.Range("L83").Value = 123500.004
$ws.Range("N83").Value = -133484.004

# Sheet 2 (ARM), row 116
$ws = $wb.Worksheets.Item(2)
$ws.Range("H116").Value = 4248.9644
$ws.Range("I116").Value = 4086.2917
$ws.Range("K116").Value = 4086.2917
$ws.Range("M116").Value = -1792.2917

# Sheet 3 (BSM), row 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("H3").Value = 4248.9644
$ws.Range("I3").Value = 4086.2917
$ws.Range("K3").Value = 4086.2917
$ws.Range("M3").Value = -3972.2917

# Sheet 3 (BSM), row 80
$ws = $wb.Worksheets.Item(3)
$ws.Range("H80").Value = 410.6
$ws.Range("I80").Value = 427.875
$ws.Range("K80").Value = 427.875
$ws.Range("M80").Value = 570.125

# Sheet 3 (BSM), row 82
$ws = $wb.Worksheets.Item(3)
$ws.Range("H82").Value = 24036
$ws.Range("J82").Value = 36245.223
$ws.Range("L82").Value = 36245.223
$ws.Range("N82").Value = -37011.223

# Sheet 3 (BSM), row 83
$ws = $wb.Worksheets.Item(3)
$ws.Range("H83").Value = 410.6
$ws.Range("I83").Value = 427.875
$ws.Range("K83").Value = 2139.375
$ws.Range("M83").Value = 2852.625

# Sheet 3 (BSM), row 85
$ws = $wb.Worksheets.Item(3)
$ws.Range("H85").Value = 24036
$ws.Range("J85").Value = 36245.223
$ws.Range("L85").Value = 36245.223
$ws.Range("N85").Value = -38897.223

# Sheet 3 (BSM), row 107
$ws = $wb.Worksheets.Item(3)
$ws.Range("H107").Value = 729.7143
$ws.Range("I107").Value = 726.3333
$ws.Range("K107").Value = 726.3333
$ws.Range("M107").Value = 1193.6667

# Sheet 3 (BSM), row 134
$ws = $wb.Worksheets.Item(3)
$ws.Range("H134").Value = 6392.257
$ws.Range("I134").Value = 5727.4365
$ws.Range("J134").Value = 8316.736999999999
$ws.Range("K134").Value = 17182.3095
$ws.Range("L134").Value = 24950.211
$ws.Range("M134").Value = -14647.3095
$ws.Range("N134").Value = -30020.211

# Sheet 4 (CRP), row 51
$ws = $wb.Worksheets.Item(4)
$ws.Range("H51").Value = 19645.166
$ws.Range("J51").Value = 19645.166
$ws.Range("L51").Value = 19645.166
$ws.Range("N51").Value = -21117.166

# Sheet 4 (CRP), row 61
$ws = $wb.Worksheets.Item(4)
$ws.Range("H61").Value = 19645.166
$ws.Range("J61").Value = 19645.166
$ws.Range("L61").Value = 19645.166
$ws.Range("N61").Value = -20341.166

# Sheet 4 (CRP), row 105
$ws = $wb.Worksheets.Item(4)
$ws.Range("H105").Value = 2056.125
$ws.Range("I105").Value = 2097.5264
$ws.Range("J105").Value = 1898.8
$ws.Range("K105").Value = 2097.5264
$ws.Range("L105").Value = 1898.8
$ws.Range("M105").Value = -350.5264000000002
$ws.Range("N105").Value = -5392.8

# Sheet 4 (CRP), row 134
$ws = $wb.Worksheets.Item(4)
$ws.Range("H134").Value = 5250.3335
$ws.Range("I134").Value = 1985
$ws.Range("J134").Value = 11781
$ws.Range("K134").Value = 5955
$ws.Range("L134").Value = 35343
$ws.Range("M134").Value = -3420
$ws.Range("N134").Value = -40413

# Sheet 5 (CUL), row 131
$ws = $wb.Worksheets.Item(5)
$ws.Range("H131").Value = 2642.8262
$ws.Range("J131").Value = 3467.5
$ws.Range("L131").Value = 10402.5
$ws.Range("N131").Value = -20482.5

# Sheet 6 (GSM), row 63
$ws = $wb.Worksheets.Item(6)
$ws.Range("H63").Value = 50000
$ws.Range("J63").Value = 50000
$ws.Range("L63").Value = 50000
$ws.Range("N63").Value = -51372

# Sheet 6 (GSM), row 66
$ws = $wb.Worksheets.Item(6)
$ws.Range("H66").Value = 50000
$ws.Range("J66").Value = 50000
$ws.Range("L66").Value = 150000
$ws.Range("N66").Value = -156864

# Sheet 6 (GSM), row 113
$ws = $wb.Worksheets.Item(6)
$ws.Range("H113").Value = 189766.17
$ws.Range("I113").Value = 226969.6
$ws.Range("J113").Value = 3749
$ws.Range("K113").Value = 226969.6
$ws.Range("L113").Value = 3749
$ws.Range("M113").Value = -224799.6
$ws.Range("N113").Value = -8089

# Sheet 7 (LTW), row 61
$ws = $wb.Worksheets.Item(7)
$ws.Range("H61").Value = 5000
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 5000
$ws.Range("N61").Value = -5404
$ws.Range("M61").Value = ""

# Sheet 7 (LTW), row 93
$ws = $wb.Worksheets.Item(7)
$ws.Range("H93").Value = 3589.8
$ws.Range("I93").Value = 3439.2632
$ws.Range("K93").Value = 3439.2632
$ws.Range("M93").Value = -2191.2632

# Sheet 7 (LTW), row 113
$ws = $wb.Worksheets.Item(7)
$ws.Range("H113").Value = 5000
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 5000
$ws.Range("N113").Value = -9340
$ws.Range("M113").Value = ""

# Sheet 8 (WVR), row 4
$ws = $wb.Worksheets.Item(8)
$ws.Range("H4").Value = 874.5
$ws.Range("J4").Value = 866
$ws.Range("L4").Value = 866
$ws.Range("N4").Value = -1092

# Sheet 8 (WVR), row 62
$ws = $wb.Worksheets.Item(8)
$ws.Range("H62").Value = 9049.166999999999
$ws.Range("J62").Value = 10124.25
$ws.Range("L62").Value = 10124.25
$ws.Range("N62").Value = -11372.25

# Sheet 8 (WVR), row 65
$ws = $wb.Worksheets.Item(8)
$ws.Range("H65").Value = 9049.166999999999
$ws.Range("J65").Value = 10124.25
$ws.Range("L65").Value = 50621.25
$ws.Range("N65").Value = -56861.25

# Sheet 8 (WVR), row 107
$ws = $wb.Worksheets.Item(8)
$ws.Range("H107").Value = 1898
$ws.Range("I107").Value = 1692.5883
$ws.Range("K107").Value = 5077.7649
$ws.Range("M107").Value = -3157.7649
